$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the reference label in D2 (was "Graham2018BIA", now "graham2018budget")
$ws.Range("D2").Value = "graham2018budget"

# Reflect the active selection moving to the edited cell
$ws.Range("D2").Select()
